$wb = $excel.ActiveWorkbook

# Per-sheet Values (column B) updates: raw counts -> normalized percentages
$updates = @{
    "D Green" = @{
        5 = 9.130434782608695
        7 = 4.201680672268908
        9 = 7.327586206896551
    }
    "Green" = @{
        2 = 3.71900826446281
        3 = 11.50793650793651
        4 = 40.08264462809917
        5 = 27.82608695652174
        6 = 22.17741935483871
        7 = 34.03361344537815
        8 = 62.82051282051282
        9 = 36.63793103448276
        10 = 21.36752136752137
        11 = 34.34782608695652
        12 = 38.30645161290323
        13 = 19.76744186046512
        14 = 27.90697674418605
        15 = 40.0
    }
    "Yellow" = @{
        2 = 42.56198347107438
        3 = 29.36507936507937
        4 = 20.24793388429752
        5 = 34.34782608695652
        6 = 33.46774193548387
        7 = 43.27731092436975
        8 = 14.52991452991453
        9 = 23.70689655172414
        10 = 34.18803418803419
        11 = 35.21739130434783
        12 = 30.24193548387097
        13 = 39.53488372093023
        14 = 25.1937984496124
        15 = 25.65217391304348
    }
    "Orange" = @{
        2 = 26.8595041322314
        3 = 22.22222222222222
        4 = 9.090909090909092
        5 = 8.26086956521739
        6 = 10.08064516129032
        7 = 3.361344537815126
        8 = 5.982905982905983
        9 = 11.63793103448276
        10 = 12.39316239316239
        11 = 11.73913043478261
        12 = 10.48387096774194
        13 = 14.34108527131783
        14 = 17.82945736434109
        15 = 14.78260869565217
    }
    "Brown" = @{
        2 = 4.132231404958678
        3 = 10.31746031746032
        4 = 4.132231404958678
        5 = 6.956521739130435
        6 = 13.70967741935484
        7 = 2.521008403361345
        8 = 4.273504273504274
        9 = 5.172413793103448
        10 = 9.82905982905983
        11 = 3.91304347826087
        12 = 4.838709677419355
        13 = 6.589147286821706
        14 = 6.976744186046512
        15 = 9.565217391304348
    }
    "Red" = @{
        2 = 16.94214876033058
        3 = 24.20634920634921
        4 = 21.48760330578512
        5 = 8.26086956521739
        6 = 17.33870967741936
        7 = 6.722689075630252
        8 = 8.974358974358974
        9 = 9.482758620689655
        10 = 17.09401709401709
        11 = 9.565217391304348
        12 = 12.09677419354839
        13 = 13.56589147286822
        14 = 18.9922480620155
        15 = 4.782608695652174
    }
    "Default Red" = @{
        2 = 5.785123966942149
        3 = 2.380952380952381
        4 = 4.958677685950413
        5 = 5.217391304347826
        6 = 3.225806451612903
        7 = 5.882352941176471
        8 = 3.418803418803419
        9 = 6.03448275862069
        10 = 5.128205128205128
        11 = 5.217391304347826
        12 = 4.032258064516129
        13 = 6.2015503875969
        14 = 3.10077519379845
        15 = 5.217391304347826
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowVals = $updates[$sheetName]
    foreach ($row in $rowVals.Keys) {
        $ws.Cells.Item([int]$row, 2).Value = $rowVals[$row]
    }
}

# "Blue" sheet: column B (Values) had no data at all; add zeros for rows 2-15
$wsBlue = $wb.Worksheets.Item("Blue")
for ($r = 2; $r -le 15; $r++) {
    $wsBlue.Cells.Item($r, 2).Value = 0
}
